# Apply the daily cryptos-list refresh (prices / 1h-volume %, and two
# re-ranked coins swapping places) coming from the GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '69.013.40'
$ws.Range('E2').Value = '  +0.41%  '

# Row 3
$ws.Range('D3').Value = '2.473.87'
$ws.Range('E3').Value = '  +0.57%  '

# Row 4
$ws.Range('E4').Value = '  +0.02%  '

# Row 5
$ws.Range('D5').Value = "'561.15"  # keep as text (avoid Excel auto-numbering)
$ws.Range('E5').Value = '  -0.66%  '

# Row 6
$ws.Range('D6').Value = "'162.12"  # keep as text (avoid Excel auto-numbering)
$ws.Range('E6').Value = '  -0.81%  '

# Row 7
$ws.Range('E7').Value = '  +0.07%  '

# Row 8
$ws.Range('D8').Value = "'0.506"  # keep as text (avoid Excel auto-numbering)
$ws.Range('E8').Value = '  -0.20%  '

# Row 9
$ws.Range('D9').Value = "'0.152"  # keep as text (avoid Excel auto-numbering)
$ws.Range('E9').Value = '  -0.40%  '

# Row 10
$ws.Range('E10').Value = '  +0.54%  '

# Row 11
$ws.Range('E11').Value = '  -3.01%  '

# Row 12
$ws.Range('E12').Value = '  +1.03%  '

# Row 13
$ws.Range('E13').Value = '  -0.12%  '

# Row 14
$ws.Range('D14').Value = '68.883.12'
$ws.Range('E14').Value = '  +0.48%  '

# Row 15
$ws.Range('E15').Value = '  -1.89%  '

# Row 16
$ws.Range('D16').Value = "'23.64"  # keep as text (avoid Excel auto-numbering)
$ws.Range('E16').Value = '  -0.44%  '

# Row 17
$ws.Range('D17').Value = '2.483.95'
$ws.Range('E17').Value = '  +1.80%  '

# Row 18
$ws.Range('D18').Value = "'10.72"  # keep as text (avoid Excel auto-numbering)
$ws.Range('E18').Value = '  -2.95%  '

# Row 19
$ws.Range('D19').Value = "'336.37"  # keep as text (avoid Excel auto-numbering)
$ws.Range('E19').Value = '  -2.62%  '

# Row 20
$ws.Range('E20').Value = '  -3.60%  '

# Row 21
$ws.Range('E21').Value = '  -1.02%  '

# Row 22
$ws.Range('E22').Value = '  +0.05%  '

# Row 23
$ws.Range('E23').Value = '  -1.12%  '

# Row 24
$ws.Range('D24').Value = "'66.83"  # keep as text (avoid Excel auto-numbering)
$ws.Range('E24').Value = '  -2.03%  '

# Row 25
$ws.Range('B25').Value = 'NEARProtocol'
$ws.Range('C25').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D25').Value = "'3.67"  # keep as text (avoid Excel auto-numbering)
$ws.Range('E25').Value = '  -2.78%  '

# Row 26
$ws.Range('B26').Value = 'WrappedeETH'
$ws.Range('C26').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D26').Value = '2.601.08'
$ws.Range('E26').Value = '  +0.70%  '

# Row 27
$ws.Range('E27').Value = '  -4.31%  '

# Row 28
$ws.Range('D28').Value = "'8.21"  # keep as text (avoid Excel auto-numbering)
$ws.Range('E28').Value = '  -0.78%  '

# Row 29
$ws.Range('D29').Value = '0.0₃0819'
$ws.Range('E29').Value = '  -3.18%  '

# Row 30
$ws.Range('D30').Value = "'7.21"  # keep as text (avoid Excel auto-numbering)
$ws.Range('E30').Value = '  -1.75%  '

# Row 31
$ws.Range('D31').Value = "'0.999"  # keep as text (avoid Excel auto-numbering)
$ws.Range('E31').Value = '  -0.11%  '

# Row 32
$ws.Range('D32').Value = "'430.16"  # keep as text (avoid Excel auto-numbering)
$ws.Range('E32').Value = '  -1.75%  '

# Row 33
$ws.Range('E33').Value = '  -4.41%  '

# Row 34
$ws.Range('E34').Value = '  -4.33%  '

# Row 35
$ws.Range('D35').Value = "'158.36"  # keep as text (avoid Excel auto-numbering)
$ws.Range('E35').Value = '  +0.87%  '

# Row 36
$ws.Range('D36').Value = "'19.02"  # keep as text (avoid Excel auto-numbering)
$ws.Range('E36').Value = '  +0.01%  '

# Row 37
$ws.Range('E37').Value = '  -0.08%  '

# Row 38
$ws.Range('D38').Value = "'0.110"  # keep as text (avoid Excel auto-numbering)
$ws.Range('E38').Value = '  -0.15%  '

# Row 39
$ws.Range('E39').Value = '  -0.98%  '

# Row 40
$ws.Range('E40').Value = '  -2.52%  '

# Row 41
$ws.Range('E41').Value = '  -2.31%  '

# Row 42
$ws.Range('E42').Value = '  -4.60%  '

# Row 43
$ws.Range('D43').Value = "'1.08"  # keep as text (avoid Excel auto-numbering)
$ws.Range('E43').Value = '  -2.63%  '

# Row 44
$ws.Range('E44').Value = '  -1.98%  '

# Row 45
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').Value = "'3.35"  # keep as text (avoid Excel auto-numbering)
$ws.Range('E45').Value = '  -1.01%  '

# Row 46
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = "'131.61"  # keep as text (avoid Excel auto-numbering)
$ws.Range('E46').Value = '  -3.13%  '

# Row 47
$ws.Range('E47').Value = '  -0.95%  '

# Row 48
$ws.Range('D48').Value = "'0.485"  # keep as text (avoid Excel auto-numbering)
$ws.Range('E48').Value = '  -1.32%  '

# Row 49
$ws.Range('E49').Value = '  -0.48%  '

# Row 50
$ws.Range('E50').Value = '  -0.47%  '

# Row 51
$ws.Range('E51').Value = '  +0.06%  '
